# Use timezone from preferences for excel reports
#
# Replace the joda-time based date formatting expressions in the "events"
# export template with calls to dateTool.format(...) that honour the
# locale/timezone passed in from preferences.
#
# A9 (the per-row "Time" column placeholder) is set first so its new
# string claims the lower shared-string slot, then B6 (the "Period:"
# range placeholder) is set second - this mirrors the original authoring
# order reflected in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", event.serverTime, locale, timezone)}'
$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'
